# Atualizei o video de finalizacao e os testes
# Add per-element description texts to the "description" column (E) of
# Planilha1 for the rows covering the Proton-Proton chain / Triple-Alpha /
# CNO-I products that were still missing a description.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rows 7-13 need the same "description" cell style already used by row 6
# (font used for the description column). Clone it onto the new cells
# before writing their text so the whole E6:E13 block shares formatting.
$ws1.Range("E6").Copy()
$ws1.Range("E7:E13").PasteSpecial(-4122)

# --- New description strings -------------------------------------------------
$ws1.Range("E3").Value  = "O Hélio-3 é um isótopo leve do hélio essencial em reações nucleares como a Cadeia Próton-Próton."
$ws1.Range("E5").Value  = "O Berílio-7 é um isótopo radioativo formado na Cadeia Próton-Próton."
$ws1.Range("E6").Value  = "O Lítio-7 é um isótopo estável formado na nucleossíntese estelar."
$ws1.Range("E7").Value  = "O Hélio-4 é um núcleo estável formado pela fusão de prótons e nêutrons."
$ws1.Range("E8").Value  = "O Boro-8 é um isótopo radioativo produzido em reações nucleares estelares."
$ws1.Range("E9").Value  = "O Berílio-8 é um isótopo extremamente instável formado em reações nucleares estelares."
$ws1.Range("E10").Value = "O Berílio-8 é um isótopo extremamente instável formado em reações nucleares estelares."
$ws1.Range("E11").Value = "O Carbono-12 é um isótopo estável formado no ciclo triple-alfa em estrelas e fundamental para a formação da vida."
$ws1.Range("E13").Value = "O Carbono-13 é um isótopo estável do carbono que surge do decaimento do Nitrogênio-13."
$ws1.Range("E12").Value = "O Nitrogênio-13 é um isótopo radioativo do nitrogênio."

# Leave the selection where the author finished editing.
$ws1.Range("E13").Select()
